# Stratos_Test_Data.xlsx - "Updated scripts for resolving conflicts"
#
# Applies the cell / view / column-width edits described by the target
# OOXML diff to the workbook that is already open as $excel.ActiveWorkbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "SitesAndCheckpoints": data edits + column width
# ---------------------------------------------------------------------
$wsSites = $wb.Worksheets.Item("SitesAndCheckpoints")
$wsSites.Range("A2").Value = "VatikaBusiness"
$wsSites.Range("B2").Value = "gate1"
$wsSites.Range("G2").Value = "VatikaBusiness@gmail.com"
$wsSites.Range("L2").Value = "VCentre01@gmail.com"
$wsSites.Columns.Item(12).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Sheet "ClientModule": data edits
# ---------------------------------------------------------------------
$wsClient = $wb.Worksheets.Item("ClientModule")
$wsClient.Range("A2").Value = "ola"
$wsClient.Range("C2").Value = "Henry@gmail.com"
$wsClient.Range("A3").Value = "Henry"
$wsClient.Range("C3").Value = "Henry@gmail.com"

# ---------------------------------------------------------------------
# Sheet "EquipmentModule": data edits, new column H, widened columns
# ---------------------------------------------------------------------
$wsEquip = $wb.Worksheets.Item("EquipmentModule")
$wsEquip.Range("A2").Value = "cocacola36"
$wsEquip.Range("B2").Value = "KUKKI-193"
$wsEquip.Range("C2").Value = 453013318
$wsEquip.Range("D2").Value = 1569901154
$wsEquip.Range("F2").Value = 5419924664
$wsEquip.Range("B3").Value = "harry-105"
$wsEquip.Range("A4").Value = "Dietcoke26"
$wsEquip.Range("B4").Value = "Diet-24"
$wsEquip.Range("C4").Value = 1611199286
$wsEquip.Range("D4").Value = 191991678
$wsEquip.Range("F4").Value = 181199443
$wsEquip.Range("G4").Value = "hbfdfhfhjhfjhywwwyrhdddkdrhff"
$wsEquip.Range("H4").Value = "hbfdfhfhjhfjhywwwyrhdddkdrhff"
$wsEquip.Columns.Item(7).ColumnWidth = 50
$wsEquip.Columns.Item(8).ColumnWidth = 40.666666666666664

# ---------------------------------------------------------------------
# View state: selections on every sheet, restore Registration's scroll
# position, and finish by activating SitesAndCheckpoints (so it becomes
# the sheet that is both selected and the one left on screen).
# ---------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("LOGIN")
$wsLogin.Range("C6").Select()

$wsReg = $wb.Worksheets.Item("Registration")
$wsReg.Range("A1").Select()
$wsReg.Range("A8").Select()

$wsClient.Range("L10").Select()

$wsEquip.Range("A12").Select()

$wsSites.Range("B9").Select()
